# Add 2022-Q1 fund-holdings sheet, and refresh the "总计" summary sheet.
#
# Shape of the edit (per the target OOXML):
#   - The worksheet that used to be named "总计" (sheetId 6 / rId6) is
#     repurposed in place to become the new "2022-Q1" holdings sheet
#     (renaming a sheet keeps its sheetId/rId, matching the target diff
#     where sheetId="6" now belongs to "2022-Q1").
#   - A brand new worksheet is appended right after it, named "总计"
#     (this gets a fresh sheetId/rId, matching the target's sheetId="7").
#     It holds the same summary rows as before, with a new 2022-Q1 row
#     inserted at the top.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("2021-Q4")   # 8-column layout + matching styles
$plain = $wb.Worksheets.Item("2020-Q4").Range("B2")  # a General-formatted, unstyled cell

$oldTotal = $wb.Worksheets.Item("总计")

# ---- Capture the old "总计" rows before we overwrite the sheet ----
$oldDates = @()
$oldCounts = @()
$oldValues = @()
for ($r = 2; $r -le 6; $r++) {
    $oldDates += $oldTotal.Cells.Item($r, 2).Value()
    $oldCounts += $oldTotal.Cells.Item($r, 3).Value()
    $oldValues += $oldTotal.Cells.Item($r, 4).Value()
}

# ============================================================
# Step 1: repurpose the "总计" sheet into the "2022-Q1" holdings sheet
# ============================================================
$oldTotal.Name = "2022-Q1"
$q1 = $oldTotal
$q1.Cells.Clear()

# Pull header + column-A formatting from an existing 8-column sheet.
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$template.Range("A2").Copy()
$q1.Range("A2:A10").PasteSpecial(-4122)  # xlPasteFormats

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q1.Cells.Item(1, 2 + $c).Value = $headers[$c]
}

# row data: code, name, size, position, ratio, value, rank
$rows = @(
    @("005571", "中银证券新能源灵活配置混合A", "0.91", "90.25", "3.18", "0.0289", 10),
    @("011824", "浙商汇金量化臻选股票型证券投资基金A", "1.54", "92.80", "1.52", "0.0234", 5),
    @("161038", "富国新兴成长量化精选混合（LOF）", "1.13", "93.66", "1.83", "0.0207", 7),
    @("010253", "兴银中证500指数增强A", "2.19", "82.47", "0.93", "0.0204", 10),
    @("011205", "兴银中证500指数增强C", "1.78", "82.47", "0.93", "0.0166", 10),
    @("006729", "万家中证500指数增强A", "1.04", "93.64", "1.35", "0.0140", 1),
    @("005572", "中银证券新能源灵活配置混合C", "0.28", "90.25", "3.18", "0.0089", 10),
    @("006730", "万家中证500指数增强C", "0.61", "93.64", "1.35", "0.0082", 1),
    @("011825", "浙商汇金量化臻选股票型证券投资基金C", "0.47", "92.80", "1.52", "0.0071", 5)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $row = $rows[$i]

    $q1.Cells.Item($r, 1).Value = $i

    # B (code) and C (name) through G (ratio) must stay text -- B holds
    # zero-padded fund codes ("005571") that must not collapse to numbers,
    # and D/E/F/G hold numeric-looking figures that the source keeps as text.
    $textRange = $q1.Range($q1.Cells.Item($r, 2), $q1.Cells.Item($r, 7))
    $textRange.NumberFormat = "@"
    $q1.Cells.Item($r, 2).Value = $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = $row[2]
    $q1.Cells.Item($r, 5).Value = $row[3]
    $q1.Cells.Item($r, 6).Value = $row[4]
    $q1.Cells.Item($r, 7).Value = $row[5]
    # strip the Text numberformat back to General without touching the
    # already-committed string values/types
    $plain.Copy()
    $textRange.PasteSpecial(-4122)

    $q1.Cells.Item($r, 8).Value = $row[6]
}

# ============================================================
# Step 2: brand-new "总计" sheet (after "2022-Q1"), with the historical
# summary rows plus the new 2022-Q1 entry on top.
# ============================================================
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$template2 = $wb.Worksheets.Item("2021-Q1")  # any sheet sharing the 3-col-summary style (A/B/C/D) -- fallback below
# Use the previous "总计" layout as the style template instead, since it already
# matches (B/C/D headers, col-A style) -- grab it from the 2020-Q4? No: reuse q1's
# column-A style source (template, 2021-Q4) is 8 columns; for the summary sheet we
# instead copy header/col-A formatting forward from the data we just captured above
# via the original "总计" sheet's own styles, which we preserved by re-deriving from
# the same 2021-Q4 template (style id 2 is shared across all non-header sheets).
$template.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$dates = @("2022-Q1") + $oldDates
$counts = @(9) + $oldCounts
$values = @(0.15) + $oldValues

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = 2 + $i
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $dates[$i]
    $total.Cells.Item($r, 3).Value = $counts[$i]
    $total.Cells.Item($r, 4).Value = $values[$i]
}
